$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update product variation code string and quantity
$ws.Range("B2").Value = "PV 001"
$ws.Range("C2").Value = 1000

# Row 3: update quantity
$ws.Range("C3").Value = 1000

# Row 4 stays the same (C4 = 500)

# Add new row 5, duplicating row 4's content
$ws.Range("A5").Value = "P REV 02"
$ws.Range("B5").Value = "PV 005"
$ws.Range("C5").Value = 500

# Update the active cell selection
$ws.Range("E9").Select()
